# Reduce the results table to only the "relevant" models (Logistic Regression
# and LightGBM), updating their metrics and dropping the rest
# (Lasso, Support Vector Classifier, CART, Random Forest, XGBoost).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 stays "Logistic Regression" but with refreshed metric values
$ws.Range("B2").Value = 0.7239669421487602
$ws.Range("C2").Value = 0.7230019698031925
$ws.Range("D2").Value = 0.7239669421487602
$ws.Range("E2").Value = 0.7143325515584248

# Row 3 (previously "Lasso") becomes the "LightGBM" row, carrying over the
# metric values that used to live on row 7 for LightGBM
$ws.Range("A3").Value = "LightGBM"
$ws.Range("B3").Value = 0.7606060606060606
$ws.Range("C3").Value = 0.7597953509291424
$ws.Range("D3").Value = 0.7606060606060606
$ws.Range("E3").Value = 0.7577348917220179

# Drop the rows for models that are no longer tracked: Support Vector
# Classifier, CART, Random Forest, the old LightGBM row, and XGBoost
$ws.Rows("4:8").Delete() | Out-Null

Write-Host "Final used range:" $ws.UsedRange.Address()
